$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "28.554.55"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.565.35"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.63"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.23"
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.20"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.789.16"
$ws.Range("D14").Value = "1.560.10"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "28.558.82"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.05"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.48"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.31"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.85"
$ws.Range("E23").Value = "  -7.05%  "
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.59"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "1.395.31"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.53"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.59"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.535"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.50"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.57"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "1.701.61"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.12"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("E51").Value = "  -0.99%  "
